# Generate Report for Handback
# 0434c169-f9ee-4296-b08c-211985e4808d.md has been handed back (in sync with en-US).
# Re-sort the status rows so the just-handed-back file moves to the top of each
# sheet, and refresh its handoff/handback/target metadata + hyperlinks.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

# Row 2 -> 0434c169 (now handed back)
$ov.Cells.Item(2,1).Value = "0434c169-f9ee-4296-b08c-211985e4808d.md"
$ov.Cells.Item(2,2).Value = "e2e\0434c169-f9ee-4296-b08c-211985e4808d.md"
$ov.Cells.Item(2,7).Value = "2016-08-19 01:01:19"

# Row 3 -> ffffbc0b
$ov.Cells.Item(3,1).Value = "ffffbc0b6be1-e5b1-4324-8fd4-6296ce27873b.md"
$ov.Cells.Item(3,2).Value = "e2e\ffffbc0b6be1-e5b1-4324-8fd4-6296ce27873b.md"

# Row 4 -> ffffff8b
$ov.Cells.Item(4,1).Value = "ffffff8b4245ff-a8d3-45b9-beb0-7e3215eeb407.md"
$ov.Cells.Item(4,2).Value = "e2e\ffffff8b4245ff-a8d3-45b9-beb0-7e3215eeb407.md"
$ov.Cells.Item(4,5).Value = "Handed back: in sync with en-US"
$ov.Cells.Item(4,6).Value = "Handed back: in sync with en-US"
$ov.Cells.Item(4,7).Value = "2016-08-19 01:00:17"

$ov.Hyperlinks.Delete()
$ov.Hyperlinks.Add($ov.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/04acf0b22012a600332b7a25711a24fc32ba3f26/e2e/0434c169-f9ee-4296-b08c-211985e4808d.md", "", "", "e2e\0434c169-f9ee-4296-b08c-211985e4808d.md")
$ov.Hyperlinks.Add($ov.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/25ddfdab2bb48d46b7ef5852d591145e567e9598/e2e/ffffbc0b6be1-e5b1-4324-8fd4-6296ce27873b.md", "", "", "e2e\ffffbc0b6be1-e5b1-4324-8fd4-6296ce27873b.md")
$ov.Hyperlinks.Add($ov.Range("B4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/04acf0b22012a600332b7a25711a24fc32ba3f26/e2e/ffffff8b4245ff-a8d3-45b9-beb0-7e3215eeb407.md", "", "", "e2e\ffffff8b4245ff-a8d3-45b9-beb0-7e3215eeb407.md")

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

# Row 2 -> 0434c169 (now handed back: full target+handback info)
$zh.Cells.Item(2,1).Value = "0434c169-f9ee-4296-b08c-211985e4808d.md"
$zh.Cells.Item(2,7).Value = "0434c169-f9ee-4296-b08c-211985e4808d.4d3c4c28663a3158cdedeba94fb5dabb0d505f2e.zh-cn.xlf"
$zh.Cells.Item(2,8).Value = "2016-08-19 01:01:14"
$zh.Cells.Item(2,9).Value = "0434c169-f9ee-4296-b08c-211985e4808d.md"
$zh.Cells.Item(2,10).Value = "0434c169-f9ee-4296-b08c-211985e4808d.4d3c4c28663a3158cdedeba94fb5dabb0d505f2e.zh-cn.xlf"
$zh.Cells.Item(2,11).Value = "2016-08-19 01:01:37"

# Row 3 -> ffffbc0b
$zh.Cells.Item(3,1).Value = "ffffbc0b6be1-e5b1-4324-8fd4-6296ce27873b.md"
$zh.Cells.Item(3,6).Value = "False"

# Row 4 -> ffffff8b
$zh.Cells.Item(4,1).Value = "ffffff8b4245ff-a8d3-45b9-beb0-7e3215eeb407.md"
$zh.Cells.Item(4,3).Value = "Handed back: in sync with en-US"
$zh.Cells.Item(4,6).Value = "True"
$zh.Cells.Item(4,7).Value = "b4eb5bdb-0a74-4a67-85cd-bcf96104d89b.cf0a9dc5466e3a6b28a7dbefc032e90daf1df6d7.zh-cn.xlf"
$zh.Cells.Item(4,8).Value = "2016-08-19 01:00:03"
$zh.Cells.Item(4,9).Value = "b4eb5bdb-0a74-4a67-85cd-bcf96104d89b.md"
$zh.Cells.Item(4,10).Value = "b4eb5bdb-0a74-4a67-85cd-bcf96104d89b.cf0a9dc5466e3a6b28a7dbefc032e90daf1df6d7.zh-cn.xlf"
$zh.Cells.Item(4,11).Value = "2016-08-19 01:00:31"

$zh.Hyperlinks.Delete()
$zh.Hyperlinks.Add($zh.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/04acf0b22012a600332b7a25711a24fc32ba3f26/e2e/0434c169-f9ee-4296-b08c-211985e4808d.md", "", "", "0434c169-f9ee-4296-b08c-211985e4808d.md")
$zh.Hyperlinks.Add($zh.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/04acf0b22012a600332b7a25711a24fc32ba3f26/e2e/0434c169-f9ee-4296-b08c-211985e4808d.md", "", "", "0434c169-f9ee-4296-b08c-211985e4808d.md")
$zh.Hyperlinks.Add($zh.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/25ddfdab2bb48d46b7ef5852d591145e567e9598/e2e/ffffbc0b6be1-e5b1-4324-8fd4-6296ce27873b.md", "", "", "ffffbc0b6be1-e5b1-4324-8fd4-6296ce27873b.md")
$zh.Hyperlinks.Add($zh.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/f848dd786196444c4d3062f4df07ee9be0464b8d/e2e/b4eb5bdb-0a74-4a67-85cd-bcf96104d89b.md", "", "", "b4eb5bdb-0a74-4a67-85cd-bcf96104d89b.md")
$zh.Hyperlinks.Add($zh.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/04acf0b22012a600332b7a25711a24fc32ba3f26/e2e/ffffff8b4245ff-a8d3-45b9-beb0-7e3215eeb407.md", "", "", "ffffff8b4245ff-a8d3-45b9-beb0-7e3215eeb407.md")
$zh.Hyperlinks.Add($zh.Range("I4"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/f848dd786196444c4d3062f4df07ee9be0464b8d/e2e/b4eb5bdb-0a74-4a67-85cd-bcf96104d89b.md", "", "", "b4eb5bdb-0a74-4a67-85cd-bcf96104d89b.md")

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

# Row 2 -> 0434c169 (now handed back: full target+handback info)
$de.Cells.Item(2,1).Value = "0434c169-f9ee-4296-b08c-211985e4808d.md"
$de.Cells.Item(2,7).Value = "0434c169-f9ee-4296-b08c-211985e4808d.4d3c4c28663a3158cdedeba94fb5dabb0d505f2e.de-de.xlf"
$de.Cells.Item(2,8).Value = "2016-08-19 01:01:19"
$de.Cells.Item(2,9).Value = "0434c169-f9ee-4296-b08c-211985e4808d.md"
$de.Cells.Item(2,10).Value = "0434c169-f9ee-4296-b08c-211985e4808d.4d3c4c28663a3158cdedeba94fb5dabb0d505f2e.de-de.xlf"
$de.Cells.Item(2,11).Value = "2016-08-19 01:01:44"

# Row 3 -> ffffbc0b
$de.Cells.Item(3,1).Value = "ffffbc0b6be1-e5b1-4324-8fd4-6296ce27873b.md"
$de.Cells.Item(3,6).Value = "False"

# Row 4 -> ffffff8b
$de.Cells.Item(4,1).Value = "ffffff8b4245ff-a8d3-45b9-beb0-7e3215eeb407.md"
$de.Cells.Item(4,3).Value = "Handed back: in sync with en-US"
$de.Cells.Item(4,6).Value = "True"
$de.Cells.Item(4,7).Value = "b4eb5bdb-0a74-4a67-85cd-bcf96104d89b.cf0a9dc5466e3a6b28a7dbefc032e90daf1df6d7.de-de.xlf"
$de.Cells.Item(4,8).Value = "2016-08-19 01:00:17"
$de.Cells.Item(4,9).Value = "b4eb5bdb-0a74-4a67-85cd-bcf96104d89b.md"
$de.Cells.Item(4,10).Value = "b4eb5bdb-0a74-4a67-85cd-bcf96104d89b.cf0a9dc5466e3a6b28a7dbefc032e90daf1df6d7.de-de.xlf"
$de.Cells.Item(4,11).Value = "2016-08-19 01:00:39"

$de.Hyperlinks.Delete()
$de.Hyperlinks.Add($de.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/04acf0b22012a600332b7a25711a24fc32ba3f26/e2e/0434c169-f9ee-4296-b08c-211985e4808d.md", "", "", "0434c169-f9ee-4296-b08c-211985e4808d.md")
$de.Hyperlinks.Add($de.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/04acf0b22012a600332b7a25711a24fc32ba3f26/e2e/0434c169-f9ee-4296-b08c-211985e4808d.md", "", "", "0434c169-f9ee-4296-b08c-211985e4808d.md")
$de.Hyperlinks.Add($de.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/25ddfdab2bb48d46b7ef5852d591145e567e9598/e2e/ffffbc0b6be1-e5b1-4324-8fd4-6296ce27873b.md", "", "", "ffffbc0b6be1-e5b1-4324-8fd4-6296ce27873b.md")
$de.Hyperlinks.Add($de.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/bd41dd78b80810e57c0e1a72aa0f221b38f63e45/e2e/b4eb5bdb-0a74-4a67-85cd-bcf96104d89b.md", "", "", "b4eb5bdb-0a74-4a67-85cd-bcf96104d89b.md")
$de.Hyperlinks.Add($de.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/04acf0b22012a600332b7a25711a24fc32ba3f26/e2e/ffffff8b4245ff-a8d3-45b9-beb0-7e3215eeb407.md", "", "", "ffffff8b4245ff-a8d3-45b9-beb0-7e3215eeb407.md")
$de.Hyperlinks.Add($de.Range("I4"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/bd41dd78b80810e57c0e1a72aa0f221b38f63e45/e2e/b4eb5bdb-0a74-4a67-85cd-bcf96104d89b.md", "", "", "b4eb5bdb-0a74-4a67-85cd-bcf96104d89b.md")
